$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New region names (column A), rows 2-9
$regions = @(
    "Rio Grande do Norte",
    "Distrito Federal",
    "Sergipe",
    "Pernambuco",
    "Maranhão",
    "Amapá",
    "Nordeste",
    "Brasil"
)

# New values (column C), rows 2-9
$values = @(
    2.519029930879796,
    2.46714858044065,
    1.753979975905736,
    1.693076222083292,
    1.472266571610874,
    1.352069305749822,
    1.129499231758544,
    0.8701736886411879
)

# New "Colocação" (column D) - only rows 2-7 have a value, rows 8-9 have none
$rank = @("1º", "2º", "3º", "4º", "5º", "6º")

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $regions[$i]
    $ws.Cells.Item($row, 2).Value = "Diferença 2024/01 - 2023/01"
    $ws.Cells.Item($row, 3).Value = $values[$i]
    if ($i -lt 6) {
        $ws.Cells.Item($row, 4).Value = $rank[$i]
    }
}

# Remove the old D8 value ("11º") since new row 8 (Nordeste) has no Colocação value,
# and old row D9/D10 were already empty.
$ws.Cells.Item(8, 4).ClearContents()

# Remove the old row 10 (Brasil) entirely - data now ends at row 9
$ws.Rows.Item(10).Delete()
